# This edit inserts two new weekly price records (rows) for Coliflor at
# "Terminal Hortofrutícola Agro Chillán" right before the existing row 150,
# pushing the former rows 150-261 down to become rows 152-263 and extending
# the sheet's used range from A1:R261 to A1:R263.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 150; everything below (old rows 150-261)
# shifts down automatically to rows 152-263.
$ws.Rows("150:151").Insert()

# --- Columns that are constant across the whole dataset: set them in bulk
# for both new rows at once. ---
$ws.Range("A150:A151").Value = 7
$ws.Range("B150:B151").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C150:C151").Value = "Ñuble"
$ws.Range("E150:E151").Value = 16
$ws.Range("F150:F151").Value = 100112008
$ws.Range("G150:G151").Value = "Coliflor"
$ws.Range("H150:H151").Value = "Sin especificar"
$ws.Range("N150:N151").Value = "`$/unidad"
$ws.Range("Q150:Q151").Value = 1
$ws.Range("R150:R151").Value = "Hortaliza"

# --- New row 150 (Primera quality) ---
$ws.Range("D150").Value = 44669
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 160
$ws.Range("K150").Value = 1200
$ws.Range("L150").Value = 1300
$ws.Range("M150").Value = 1250
$ws.Range("O150").Value = "Provincia de Diguillín"
$ws.Range("P150").Value = 1250

# --- New row 151 (Segunda quality) ---
$ws.Range("D151").Value = 44669
$ws.Range("I151").Value = "Segunda"
$ws.Range("J151").Value = 60
$ws.Range("K151").Value = 1000
$ws.Range("L151").Value = 1000
$ws.Range("M151").Value = 1000
$ws.Range("O151").Value = "Provincia de Diguillín"
$ws.Range("P151").Value = 1000
